$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows (31 and 32) with new Mac-Address / machine entries,
# following the same pattern as the existing rows (regcntr_id=10001,
# incrementing machine_id, lang_code="eng", is_active=TRUE,
# cr_by="superadmin", cr_dtimes="now()")

$ws.Cells.Item(31, 1).Value = 10001
$ws.Cells.Item(31, 2).Value = 10030
$ws.Cells.Item(31, 3).Value = "eng"
$ws.Cells.Item(31, 4).Value = $true
$ws.Cells.Item(31, 5).Value = "superadmin"
$ws.Cells.Item(31, 6).Value = "now()"

$ws.Cells.Item(32, 1).Value = 10001
$ws.Cells.Item(32, 2).Value = 10031
$ws.Cells.Item(32, 3).Value = "eng"
$ws.Cells.Item(32, 4).Value = $true
$ws.Cells.Item(32, 5).Value = "superadmin"
$ws.Cells.Item(32, 6).Value = "now()"

# Move the selection to where the user last clicked, matching the
# post-edit view state (cell E31).
$ws.Range("E31").Select()
